# Auto-generated Excel COM-interop script
# Applies: float_format='%.2f' rounding to the 'Intersects' sheet column B
# (nzd0014-0000 distances), and propagates the resulting change to the
# downstream regression statistics on 'Transects' and the recomputed
# intersection coordinates on 'Intersect points'.

$wb = $excel.ActiveWorkbook
$wsIntersects = $wb.Worksheets.Item("Intersects")
$wsTransects = $wb.Worksheets.Item("Transects")
$wsPoints = $wb.Worksheets.Item("Intersect points")

# New rounded distance values (column B, "nzd0014-0000") keyed by row number
$distanceUpdates = @{
    2 = 372.71
    3 = 381.9
    4 = 373.45
    5 = 386.35
    6 = 366.59
    7 = 374.31
    8 = 376.85
    9 = 382.95
    10 = 376.27
    11 = 371.35
    12 = 359.12
    14 = 376.44
    15 = 363.09
    16 = 358.49
    17 = 360.06
    18 = 357.43
    19 = 379.24
    20 = 376.47
    21 = 364.66
    22 = 373.77
    24 = 372.89
    25 = 382.25
    26 = 378.55
    27 = 375.43
    29 = 376.84
    30 = 347.69
    31 = 375.59
    32 = 357.02
    34 = 370.4
    35 = 370.75
    36 = 389.38
    37 = 383.86
    38 = 382.5
    39 = 370.12
    40 = 383.59
    41 = 388.2
    42 = 372.94
    43 = 364.39
    44 = 388.02
    46 = 374.53
    47 = 367.54
    48 = 365.86
    49 = 383.8
    51 = 379.66
    52 = 372.7
    53 = 387.3
    54 = 384.01
    55 = 365.4
    56 = 377.54
    57 = 360.6
    58 = 381.12
    59 = 374.05
    60 = 368.73
    61 = 372.62
    62 = 362.27
    63 = 356.67
    65 = 376.16
    66 = 372.21
    68 = 360.25
    69 = 360.09
    70 = 375.18
    71 = 377.84
    72 = 375.03
    74 = 365.6
    75 = 353.66
    78 = 374.33
    80 = 375.76
    81 = 400.13
    82 = 350.23
    83 = 376
    84 = 366.07
    85 = 361.19
    86 = 372.77
    87 = 361.06
    88 = 364.52
    89 = 373.39
    90 = 382.57
    91 = 374.59
    92 = 378.34
    93 = 369.28
    94 = 359.52
    95 = 356.62
    96 = 370.32
    98 = 366.63
    99 = 363.12
    101 = 383.46
    102 = 375.67
    104 = 377.56
    105 = 368.04
    106 = 360.53
    107 = 378.35
    108 = 380.19
    109 = 380.53
    112 = 378.56
    113 = 378.33
    114 = 376.07
    115 = 383.1
    117 = 383.43
    118 = 384.43
    119 = 383.97
    120 = 360.98
    121 = 380.25
    122 = 366.82
    123 = 378.66
    124 = 372.29
    125 = 383.61
    126 = 371.22
    127 = 373.59
    130 = 371.08
    132 = 355.75
    133 = 373.67
    134 = 371.24
    135 = 364.67
    136 = 377.24
    137 = 370.22
    138 = 378.45
    140 = 371.9
    141 = 374.32
    143 = 367.72
    144 = 371
    145 = 373.38
    146 = 367.21
    147 = 370.5
    148 = 378.22
    149 = 375.37
    150 = 370.62
    151 = 368.38
    152 = 361.42
    153 = 359.08
    154 = 364.49
    156 = 371.29
    157 = 385.29
    158 = 377.72
    159 = 381.74
    160 = 378.97
    161 = 376.31
    163 = 357.24
    164 = 364.73
    165 = 357.95
    168 = 372.52
    169 = 377.48
    170 = 371.96
    171 = 358.5
    172 = 375.91
    173 = 384.21
    174 = 370.97
    175 = 375.72
    176 = 365.74
    178 = 373.44
    179 = 369.22
    180 = 372.56
    181 = 372.1
    182 = 371.6
    183 = 373.02
    184 = 369.2
    185 = 384.99
    186 = 378.78
    187 = 380.61
    188 = 381.43
    190 = 376.45
    192 = 372.97
    193 = 380.98
    194 = 383.47
    195 = 363
    198 = 381.98
    199 = 381.53
    200 = 360.15
    201 = 378.51
    202 = 379.22
    203 = 380.74
    204 = 365.95
    205 = 381.97
    206 = 376
    207 = 376.7
    208 = 381.43
    209 = 366.24
    211 = 353.07
    212 = 371.52
    213 = 373.4
    214 = 343.27
    215 = 353.82
    216 = 369.73
    218 = 383.83
    219 = 374.07
    220 = 375.98
    221 = 381.53
}

foreach ($row in $distanceUpdates.Keys) {
    $wsIntersects.Cells.Item($row, 2).Value = $distanceUpdates[$row]
}

# Recomputed "lat,long" strings (column B) for the same rows, reflecting the
# rounded distances above
$coordUpdates = @{
    2 = "-34.72455748374235,173.09249979944377"
    3 = "-34.72453883384497,173.0925975999973"
    4 = "-34.7245559820125,173.09250767457158"
    5 = "-34.7245298031265,173.09264495715914"
    6 = "-34.72456990343459,173.09243466999737"
    7 = "-34.72455423675824,173.09251682674682"
    8 = "-34.72454908216608,173.09254385758763"
    9 = "-34.72453670300294,173.09260877416003"
    10 = "-34.724550259199475,173.09253768519116"
    11 = "-34.72456024367697,173.09248532623516"
    12 = "-34.724585062717665,173.09235517373483"
    14 = "-34.72454991420695,173.09253949434188"
    15 = "-34.724577006184134,173.0923974227562"
    16 = "-34.724586341209054,173.09234846922948"
    17 = "-34.72458315512668,173.09236517728223"
    18 = "-34.72458849232071,173.0923371886326"
    19 = "-34.7245442319735,173.09256929211583"
    20 = "-34.72454985332592,173.09253981360376"
    21 = "-34.72457382009503,173.09241413080517"
    22 = "-34.724555332615644,173.09251108003218"
    24 = "-34.72455711845676,173.09250171501546"
    25 = "-34.72453812356441,173.09260132471826"
    26 = "-34.72454563223885,173.09256194909335"
    27 = "-34.72455196386796,173.09252874585812"
    29 = "-34.72454910245976,173.09254375116697"
    30 = "-34.7246082581468,173.09223353481892"
    31 = "-34.72455163916924,173.09253044858826"
    32 = "-34.72458932435418,173.0923328253827"
    34 = "-34.72456217157147,173.0924752162727"
    35 = "-34.724561461294634,173.09247894099576"
    36 = "-34.72452365410992,173.09267720259155"
    37 = "-34.724534856272356,173.0926184584339"
    38 = "-34.72453761622107,173.0926039852332"
    39 = "-34.72456273979283,173.0924722364942"
    40 = "-34.7245354042035,173.09261558507797"
    41 = "-34.72452604877746,173.092664644965"
    42 = "-34.724557016988534,173.0925022471187"
    43 = "-34.72457436802135,173.0924112574465"
    44 = "-34.724526414065615,173.09266272939476"
    46 = "-34.72455379029773,173.09251916800088"
    47 = "-34.72456797554348,173.09244477996177"
    48 = "-34.724571384866174,173.0924269012876"
    49 = "-34.72453497803484,173.09261781991037"
    51 = "-34.72454337963785,173.09257376178155"
    52 = "-34.72455750403599,173.09249969302311"
    53 = "-34.72452787521794,173.09265506711372"
    54 = "-34.72453455186615,173.09262005474275"
    55 = "-34.72457231837073,173.0924220059361"
    56 = "-34.72454768190181,173.0925512006107"
    57 = "-34.724582059276145,173.0923709240007"
    58 = "-34.72454041675551,173.09258929919034"
    59 = "-34.724554764393325,173.09251405981016"
    60 = "-34.72456556060499,173.0924574440217"
    61 = "-34.72455766638513,173.09249884165794"
    62 = "-34.72457867025524,173.09238869625875"
    63 = "-34.724590034626544,173.0923291006571"
    65 = "-34.724550482429905,173.09253651456424"
    66 = "-34.7245584984244,173.0924944784113"
    68 = "-34.72458276954968,173.09236719927577"
    69 = "-34.72458309424609,173.09236549654437"
    70 = "-34.72455247120963,173.09252608534226"
    71 = "-34.72454707309113,173.09255439322936"
    72 = "-34.72455277561462,173.09252448903274"
    74 = "-34.724571912499194,173.0924241343498"
    75 = "-34.72459614296406,173.0922970680145"
    78 = "-34.72455419617092,173.09251703958807"
    80 = "-34.724551294176834,173.092532257739"
    81 = "-34.724501838223176,173.09279160466446"
    82 = "-34.724603103617476,173.092260565695"
    83 = "-34.72455080712869,173.09253481183418"
    84 = "-34.72457095870097,173.09242913612198"
    85 = "-34.72458086195768,173.09237720282258"
    86 = "-34.72455736198049,173.0925004379677"
    87 = "-34.724581125773646,173.09237581935335"
    88 = "-34.724574104205,173.0924126409155"
    89 = "-34.72455610377438,173.09250703604772"
    90 = "-34.724537474164926,173.09260473017738"
    91 = "-34.724553668535776,173.0925198065247"
    92 = "-34.72454605840646,173.09255971426035"
    93 = "-34.72456444445652,173.09246329715842"
    94 = "-34.724584250976925,173.0923594305636"
    95 = "-34.724590136094,173.09232856855348"
    96 = "-34.72456233392043,173.09247436490742"
    98 = "-34.72456982226025,173.09243509568012"
    99 = "-34.72457694530346,173.09239774201828"
    101 = "-34.724535668022156,173.09261420161027"
    102 = "-34.724551476819876,173.09253129995332"
    104 = "-34.72454764131445,173.09255141345193"
    105 = "-34.7245669608636,173.09245010099545"
    106 = "-34.72458220133086,173.0923701790557"
    107 = "-34.72454603811277,173.09255982068098"
    108 = "-34.72454230407121,173.09257940207394"
    109 = "-34.72454161408492,173.09258302037458"
    112 = "-34.72454561194515,173.09256205551395"
    113 = "-34.72454607870017,173.09255960783975"
    114 = "-34.724550665072975,173.09253555677859"
    115 = "-34.72453639859686,173.09261037046895"
    117 = "-34.72453572890338,173.09261388234847"
    118 = "-34.72453369952864,173.09262452440746"
    119 = "-34.724534633041145,173.0926196290604"
    120 = "-34.72458128812192,173.0923749679877"
    121 = "-34.72454218230892,173.0925800405976"
    122 = "-34.72456943668208,173.09243711767303"
    123 = "-34.724545409008165,173.09256311972013"
    124 = "-34.72455833607528,173.0924953297765"
    125 = "-34.72453536361599,173.09261579791914"
    126 = "-34.724560507494154,173.09248394276665"
    127 = "-34.72455569790139,173.0925091644606"
    130 = "-34.724560791604965,173.09248245287748"
    132 = "-34.7245919016276,173.0923193099496"
    133 = "-34.72455553555216,173.09251001582575"
    134 = "-34.7245604669069,173.09248415560793"
    135 = "-34.72457379980147,173.09241423722588"
    136 = "-34.72454829071241,173.09254800799198"
    137 = "-34.72456253685663,173.0924733007008"
    138 = "-34.72454583517582,173.09256088488718"
    140 = "-34.72455912752715,173.0924911793711"
    141 = "-34.72455421646458,173.09251693316745"
    143 = "-34.724567610258745,173.0924466955339"
    144 = "-34.72456095395399,173.09248160151222"
    145 = "-34.72455612406804,173.0925069296271"
    146 = "-34.72456864523208,173.09244126807945"
    147 = "-34.724561968635236,173.0924762804793"
    148 = "-34.724546301930815,173.09255843721294"
    149 = "-34.72455208562996,173.09252810733432"
    150 = "-34.724561725111755,173.09247755752722"
    151 = "-34.724566270881176,173.09245371929828"
    152 = "-34.724580395206324,173.09237965049883"
    153 = "-34.72458514389174,173.09235474805197"
    154 = "-34.724574165085684,173.09241232165343"
    156 = "-34.724560365438755,173.09248468771122"
    157 = "-34.7245319542656,173.09263367657772"
    158 = "-34.72454731661543,173.0925531161819"
    159 = "-34.72453915854461,173.09259589726767"
    160 = "-34.724544779903475,173.09256641875922"
    161 = "-34.724550178024764,173.0925381108737"
    163 = "-34.72458887789722,173.09233516663878"
    164 = "-34.72457367804005,173.09241487575"
    165 = "-34.72458743705852,173.09234272251038"
    168 = "-34.72455786932155,173.09249777745146"
    169 = "-34.724547803663945,173.09255056208696"
    170 = "-34.724559005765336,173.09249181789502"
    171 = "-34.72458632091554,173.0923485756502"
    172 = "-34.72455098977176,173.0925338540485"
    173 = "-34.724534145991164,173.0926221831545"
    174 = "-34.72456101483487,173.09248128225025"
    175 = "-34.724551375351524,173.0925318320565"
    176 = "-34.72457162838911,173.09242562423938"
    178 = "-34.724556002306144,173.09250756815095"
    179 = "-34.72456456621818,173.09246265863445"
    180 = "-34.72455778814698,173.09249820313406"
    181 = "-34.72455872165441,173.09249330778417"
    182 = "-34.724559736336175,173.09248798675154"
    183 = "-34.72455685463937,173.09250309848383"
    184 = "-34.7245646068054,173.0924624457931"
    185 = "-34.72453256307836,173.09263048396025"
    186 = "-34.72454516548378,173.09256439676753"
    187 = "-34.724541451735185,173.09258387173946"
    188 = "-34.7245397876501,173.09259259822903"
    190 = "-34.724549893913284,173.0925396007625"
    192 = "-34.7245569561076,173.09250256638063"
    193 = "-34.724540700867614,173.09258780930188"
    194 = "-34.724535647728416,173.09261430803087"
    195 = "-34.72457718882612,173.0923964649699"
    198 = "-34.72453867149514,173.0925984513621"
    199 = "-34.72453958471286,173.09259366243506"
    200 = "-34.72458297248494,173.09236613506863"
    201 = "-34.72454571341364,173.09256152341086"
    202 = "-34.72454427256091,173.09256907927463"
    203 = "-34.724541187916856,173.09258525520735"
    204 = "-34.724571202223935,173.09242785907378"
    205 = "-34.724538691788865,173.0925983449415"
    206 = "-34.72455080712869,173.09253481183418"
    207 = "-34.724549386571304,173.09254226127823"
    208 = "-34.7245397876501,173.09259259822903"
    209 = "-34.72457061371007,173.09243094527355"
    211 = "-34.7245973402784,173.0922907891903"
    212 = "-34.724559898685236,173.0924871353863"
    213 = "-34.72455608348074,173.09250714246838"
    214 = "-34.724617227825235,173.09218649682964"
    215 = "-34.72459581826859,173.09229877074642"
    216 = "-34.72456353124391,173.09246808608833"
    218 = "-34.7245349171536,173.09261813917215"
    219 = "-34.724554723806,173.09251427265144"
    220 = "-34.724550847716046,173.0925345989929"
    221 = "-34.72453958471286,173.09259366243506"
}

foreach ($row in $coordUpdates.Keys) {
    $wsPoints.Cells.Item($row, 2).Value = $coordUpdates[$row]
}

# Recomputed linear-regression statistics for transect nzd0014-0000 (row 2)
$transectUpdates = @{
    "I2" = 0.01874337809980859   # trend
    "L2" = 0.0002648221685296326   # r2_score
    "M2" = 7.151166085980697   # mae
    "N2" = 82.5588205282918   # mse
    "O2" = 9.086188448865222   # rmse
    "P2" = 372.358477219704   # intercept
}

foreach ($ref in $transectUpdates.Keys) {
    $wsTransects.Range($ref).Value = $transectUpdates[$ref]
}

Write-Output "Applied tidally_corrected.csv %.2f rounding updates"
